$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows appended through 2021-09-01 (aggiornamento fino a 1/09/2021)
$ws.Range("A358").Value = 44432
$ws.Range("B358").Value = 0
$ws.Range("C358").Value = 29
$ws.Range("D358").Value = 71.99960276081235
$ws.Range("A359").Value = 44433
$ws.Range("B359").Value = 1
$ws.Range("C359").Value = 27
$ws.Range("D359").Value = 67.03411291523909
$ws.Range("A360").Value = 44434
$ws.Range("B360").Value = 3
$ws.Range("C360").Value = 26
$ws.Range("D360").Value = 64.55136799245246
$ws.Range("A361").Value = 44435
$ws.Range("B361").Value = 9
$ws.Range("C361").Value = 31
$ws.Range("D361").Value = 76.96509260638562
$ws.Range("A362").Value = 44436
$ws.Range("B362").Value = 5
$ws.Range("C362").Value = 30
$ws.Range("D362").Value = 74.48234768359899
$ws.Range("A363").Value = 44437
$ws.Range("B363").Value = 1
$ws.Range("C363").Value = 19
$ws.Range("D363").Value = 47.17215353294603
$ws.Range("A364").Value = 44438
$ws.Range("B364").Value = 8
$ws.Range("C364").Value = 27
$ws.Range("D364").Value = 67.03411291523909
$ws.Range("A365").Value = 44439
$ws.Range("B365").Value = 2
$ws.Range("C365").Value = 29
$ws.Range("D365").Value = 71.99960276081235
$ws.Range("A366").Value = 44440
$ws.Range("B366").Value = 3
$ws.Range("C366").Value = 31
$ws.Range("D366").Value = 76.96509260638562

# Match the formatting (date style, borders, alignment) of the preceding data row
$ws.Range("A357:D357").Copy()
$ws.Range("A358:D366").PasteSpecial(-4122)
